{"js": "// Replace the contents of each lattice-multiplication table cell (in row-major\n// order) with the new exercise text. Each cell holds a single paragraph whose\n// run contains five lines (problem header, two-digit split, a separator, and\n// two lattice rows) joined by line breaks. We rebuild that text with the new\n// numbers while keeping the original run formatting (font size 32) intact by\n// doing a \"Replace\" insertText on the existing paragraph.\n\nconst NEW_CELL_TEXT = [\n  [\"61 x 31\", \"  3    1\", \"  ----\", \"6|    |\", \"1|    |\"],\n  [\"98 x 79\", \"  7    9\", \"  ----\", \"9|    |\", \"8|    |\"],\n  [\"82 x 89\", \"  8    9\", \"  ----\", \"8|    |\", \"2|    |\"],\n  [\"29 x 91\", \"  9    1\", \"  ----\", \"2|    |\", \"9|    |\"],\n  [\"62 x 52\", \"  5    2\", \"  ----\", \"6|    |\", \"2|    |\"],\n  [\"28 x 48\", \"  4    8\", \"  ----\", \"2|    |\", \"8|    |\"],\n  [\"25 x 24\", \"  2    4\", \"  ----\", \"2|    |\", \"5|    |\"],\n  [\"24 x 46\", \"  4    6\", \"  ----\", \"2|    |\", \"4|    |\"],\n  [\"54 x 24\", \"  2    4\", \"  ----\", \"5|    |\", \"4|    |\"],\n  [\"98 x 44\", \"  4    4\", \"  ----\", \"9|    |\", \"8|    |\"],\n  [\"18 x 67\", \"  6    7\", \"  ----\", \"1|    |\", \"8|    |\"],\n  [\"10 x 20\", \"  2    0\", \"  ----\", \"1|    |\", \"0|    |\"],\n  [\"10 x 67\", \"  6    7\", \"  ----\", \"1|    |\", \"0|    |\"],\n  [\"21 x 27\", \"  2    7\", \"  ----\", \"2|    |\", \"1|    |\"],\n  [\"89 x 30\", \"  3    0\", \"  ----\", \"8|    |\", \"9|    |\"],\n];\n\n// Word represents a manual line break (<w:br/>) as a vertical-tab (\\v,\n// U+000B) character inside Range/Paragraph text.\nconst LINE_BREAK = \"\\u000b\";\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Load every row's cells up front.\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nlet cellIndex = 0;\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    if (cellIndex >= NEW_CELL_TEXT.length) break;\n    const lines = NEW_CELL_TEXT[cellIndex];\n    const newText = lines.join(LINE_BREAK);\n    const para = cell.body.paragraphs.getFirst();\n    para.insertText(newText, \"Replace\");\n    cellIndex++;\n  }\n}\nawait context.sync();\n", "ps1": "# Replace the contents of each lattice-multiplication table cell (in\n# row-major order) with the new exercise text. Each cell holds a single\n# paragraph/run made up of five lines (problem header, two-digit split, a\n# separator, and two lattice rows) joined by manual line breaks\n# (Chr(11), Word's vertical-tab line-break character). Assigning directly\n# to Cell.Range.Text keeps the existing run formatting (font size 32).\n\n$lb = [char]11\n\n$newCellText = @(\n    @(\"61 x 31\", \"  3    1\", \"  ----\", \"6|    |\", \"1|    |\"),\n    @(\"98 x 79\", \"  7    9\", \"  ----\", \"9|    |\", \"8|    |\"),\n    @(\"82 x 89\", \"  8    9\", \"  ----\", \"8|    |\", \"2|    |\"),\n    @(\"29 x 91\", \"  9    1\", \"  ----\", \"2|    |\", \"9|    |\"),\n    @(\"62 x 52\", \"  5    2\", \"  ----\", \"6|    |\", \"2|    |\"),\n    @(\"28 x 48\", \"  4    8\", \"  ----\", \"2|    |\", \"8|    |\"),\n    @(\"25 x 24\", \"  2    4\", \"  ----\", \"2|    |\", \"5|    |\"),\n    @(\"24 x 46\", \"  4    6\", \"  ----\", \"2|    |\", \"4|    |\"),\n    @(\"54 x 24\", \"  2    4\", \"  ----\", \"5|    |\", \"4|    |\"),\n    @(\"98 x 44\", \"  4    4\", \"  ----\", \"9|    |\", \"8|    |\"),\n    @(\"18 x 67\", \"  6    7\", \"  ----\", \"1|    |\", \"8|    |\"),\n    @(\"10 x 20\", \"  2    0\", \"  ----\", \"1|    |\", \"0|    |\"),\n    @(\"10 x 67\", \"  6    7\", \"  ----\", \"1|    |\", \"0|    |\"),\n    @(\"21 x 27\", \"  2    7\", \"  ----\", \"2|    |\", \"1|    |\"),\n    @(\"89 x 30\", \"  3    0\", \"  ----\", \"8|    |\", \"9|    |\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$cellIndex = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        if ($cellIndex -ge $newCellText.Length) { continue }\n        $lines = $newCellText[$cellIndex]\n        $newText = [string]::Join($lb, $lines)\n        $t.Cell($r, $c).Range.Text = $newText\n        $cellIndex++\n    }\n}\n"}
